$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.90949261935522
$ws.Range("C2").Value = 23.84257251712773
$ws.Range("D2").Value = 0.5808724125471856
$ws.Range("E2").Value = 1.070126019124857
$ws.Range("F2").Value = 7.077160381832315

$ws.Range("B3").Value = 12.47492109345336
$ws.Range("C3").Value = 18.27046271825304
$ws.Range("D3").Value = 0.6795077543958578
$ws.Range("E3").Value = 1.174826161918675
$ws.Range("F3").Value = 5.747076404695129

$ws.Range("B4").Value = 11.13331018869496
$ws.Range("C4").Value = 14.2105713162964
$ws.Range("D4").Value = 0.7792304716135031
$ws.Range("E4").Value = 1.36681778862738
$ws.Range("F4").Value = 4.820984854866849
